$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 666868.9
$ws.Range("I19").Value = 1481654.2
$ws.Range("J19").Value = 226.27272
$ws.Range("K19").Value = 1481654.2
$ws.Range("L19").Value = 226.27272
$ws.Range("M19").Value = -1481479.2
$ws.Range("N19").Value = -576.2727199999999

$ws.Range("H41").Value = 824.1667
$ws.Range("I41").Value = 378
$ws.Range("K41").Value = 378
$ws.Range("M41").Value = 62

$ws.Range("H55").Value = 168.27272
$ws.Range("J55").Value = 188
$ws.Range("L55").Value = 188
$ws.Range("N55").Value = -616

$ws.Range("H106").Value = 3502.5
$ws.Range("I106").Value = 3502.5
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 3502.5
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -2871.5
$ws.Range("N106").ClearContents()

$ws.Range("H137").Value = 2975.775
$ws.Range("I137").Value = 1696.5
$ws.Range("J137").Value = 10225
$ws.Range("K137").Value = 5089.5
$ws.Range("L137").Value = 30675
$ws.Range("M137").Value = -2539.5
$ws.Range("N137").Value = -35775

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7440.65
$ws.Range("I32").Value = 3955.4626
$ws.Range("J32").Value = 14516.637
$ws.Range("K32").Value = 3955.4626
$ws.Range("L32").Value = 14516.637
$ws.Range("M32").Value = -3668.4626
$ws.Range("N32").Value = -15090.637

$ws.Range("H92").Value = 30000
$ws.Range("J92").Value = 30000
$ws.Range("L92").Value = 30000
$ws.Range("N92").Value = -34992

$ws.Range("H132").Value = 2309.8293
$ws.Range("I132").Value = 1170.8334
$ws.Range("J132").Value = 5416.1816
$ws.Range("K132").Value = 3512.5002
$ws.Range("L132").Value = 16248.5448
$ws.Range("M132").Value = -982.5001999999999
$ws.Range("N132").Value = -21308.5448

$ws.Range("H139").Value = 43406.05
$ws.Range("J139").Value = 43406.05
$ws.Range("L139").Value = 43406.05
$ws.Range("N139").Value = -53686.05

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6933.84
$ws.Range("I20").Value = 1235.8334
$ws.Range("J20").Value = 12193.538
$ws.Range("K20").Value = 1235.8334
$ws.Range("L20").Value = 12193.538
$ws.Range("M20").Value = -988.8334
$ws.Range("N20").Value = -12687.538

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H134").Value = 2729.5293
$ws.Range("I134").Value = 1512
$ws.Range("J134").Value = 7721.4
$ws.Range("K134").Value = 4536
$ws.Range("L134").Value = 23164.2
$ws.Range("M134").Value = -2001
$ws.Range("N134").Value = -28234.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3778.0715
$ws.Range("I31").Value = 1543.6666
$ws.Range("J31").Value = 7800
$ws.Range("K31").Value = 1543.6666
$ws.Range("L31").Value = 7800
$ws.Range("M31").Value = -1248.6666
$ws.Range("N31").Value = -8390

$ws.Range("H34").Value = 3778.0715
$ws.Range("I34").Value = 1543.6666
$ws.Range("J34").Value = 7800
$ws.Range("K34").Value = 1543.6666
$ws.Range("L34").Value = 7800
$ws.Range("M34").Value = -1341.6666
$ws.Range("N34").Value = -8204

$ws.Range("H81").Value = 25000
$ws.Range("J81").Value = 25000
$ws.Range("L81").Value = 25000
$ws.Range("N81").Value = -26996

$ws.Range("H84").Value = 25000
$ws.Range("J84").Value = 25000
$ws.Range("L84").Value = 75000
$ws.Range("N84").Value = -84984

$ws.Range("H99").Value = 3571.9375
$ws.Range("I99").Value = 1640.4445
$ws.Range("J99").Value = 6055.2856
$ws.Range("K99").Value = 1640.4445
$ws.Range("L99").Value = 6055.2856
$ws.Range("M99").Value = -142.4445000000001
$ws.Range("N99").Value = -9051.285599999999

$ws.Range("H122").Value = 4065.6365
$ws.Range("I122").Value = 2914.6667
$ws.Range("J122").Value = 9245
$ws.Range("K122").Value = 8744.000100000001
$ws.Range("L122").Value = 27735
$ws.Range("M122").Value = -6294.000100000001
$ws.Range("N122").Value = -32635

$ws.Range("H126").Value = 3571.9375
$ws.Range("I126").Value = 1640.4445
$ws.Range("J126").Value = 6055.2856
$ws.Range("K126").Value = 4921.333500000001
$ws.Range("L126").Value = 18165.8568
$ws.Range("M126").Value = -2451.333500000001
$ws.Range("N126").Value = -23105.8568

$ws.Range("H132").Value = 2287.1633
$ws.Range("I132").Value = 1981.2162
$ws.Range("K132").Value = 5943.6486
$ws.Range("M132").Value = -3413.6486

$ws.Range("H134").Value = 3263.889
$ws.Range("I134").Value = 3558.0222
$ws.Range("J134").Value = 2528.5557
$ws.Range("K134").Value = 10674.0666
$ws.Range("L134").Value = 7585.6671
$ws.Range("M134").Value = -8139.0666
$ws.Range("N134").Value = -12655.6671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H96").Value = 4800
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 4800
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 14400
$ws.Range("N96").Value = -18518
$ws.Range("M96").ClearContents()

$ws.Range("H98").Value = 201.5
$ws.Range("I98").Value = 103
$ws.Range("K98").Value = 309
$ws.Range("M98").Value = 1189

$ws.Range("H102").Value = 4441.6665
$ws.Range("J102").Value = 6150
$ws.Range("L102").Value = 18450
$ws.Range("N102").Value = -23318

$ws.Range("H105").Value = 5000
$ws.Range("J105").Value = 5000
$ws.Range("L105").Value = 15000
$ws.Range("N105").Value = -20242

$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("N110").ClearContents()

$ws.Range("H111").Value = 5111.75
$ws.Range("I111").Value = 223.5
$ws.Range("K111").Value = 670.5
$ws.Range("M111").Value = 2396.5

$ws.Range("H114").Value = 45455844
$ws.Range("I114").Value = 76923320
$ws.Range("J114").Value = 2819.5557
$ws.Range("K114").Value = 230769960
$ws.Range("L114").Value = 8458.667099999999
$ws.Range("M114").Value = -230766706
$ws.Range("N114").Value = -14966.6671

$ws.Range("H116").Value = 3331.6
$ws.Range("I116").Value = 414.5
$ws.Range("J116").Value = 15000
$ws.Range("K116").Value = 1243.5
$ws.Range("L116").Value = 45000
$ws.Range("N116").Value = -51884
$ws.Range("M116").Value = 2198.5

$ws.Range("H117").Value = 1664.1666
$ws.Range("I117").Value = 1397
$ws.Range("J117").Value = 3000
$ws.Range("K117").Value = 4191
$ws.Range("L117").Value = 9000
$ws.Range("M117").Value = -749
$ws.Range("N117").Value = -15884

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4146.5
$ws.Range("I122").Value = 2804.2
$ws.Range("J122").Value = 5488.8
$ws.Range("K122").Value = 8412.599999999999
$ws.Range("L122").Value = 16466.4
$ws.Range("M122").Value = -5962.599999999999
$ws.Range("N122").Value = -21366.4

$ws.Range("H126").Value = 3949.261
$ws.Range("I126").Value = 2830.66
$ws.Range("J126").Value = 5280.9287
$ws.Range("K126").Value = 8491.98
$ws.Range("L126").Value = 15842.7861
$ws.Range("M126").Value = -6021.98
$ws.Range("N126").Value = -20782.7861

$ws.Range("H132").Value = 2510.3953
$ws.Range("I132").Value = 1008.1429
$ws.Range("J132").Value = 3235.6206
$ws.Range("K132").Value = 3024.4287
$ws.Range("L132").Value = 9706.861800000001
$ws.Range("M132").Value = -494.4287000000004
$ws.Range("N132").Value = -14766.8618

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6005.857
$ws.Range("I7").Value = 3684
$ws.Range("J7").Value = 8327.714
$ws.Range("K7").Value = 3684
$ws.Range("L7").Value = 8327.714
$ws.Range("M7").Value = -3572
$ws.Range("N7").Value = -8551.714

$ws.Range("H16").Value = 574
$ws.Range("I16").Value = 574
$ws.Range("K16").Value = 574
$ws.Range("M16").Value = -404

$ws.Range("H22").Value = 7144671
$ws.Range("I22").Value = 13159242
$ws.Range("J22").Value = 2368.1875
$ws.Range("K22").Value = 13159242
$ws.Range("L22").Value = 2368.1875
$ws.Range("M22").Value = -13158947
$ws.Range("N22").Value = -2958.1875

$ws.Range("H27").Value = 7144671
$ws.Range("I27").Value = 13159242
$ws.Range("J27").Value = 2368.1875
$ws.Range("K27").Value = 13159242
$ws.Range("L27").Value = 2368.1875
$ws.Range("M27").Value = -13159135
$ws.Range("N27").Value = -2582.1875

$ws.Range("H40").Value = 5117.745
$ws.Range("I40").Value = 4660.9287
$ws.Range("K40").Value = 4660.9287
$ws.Range("M40").Value = -4524.9287

$ws.Range("H46").Value = 2225.5833
$ws.Range("I46").Value = 2579.8
$ws.Range("J46").Value = 2132.3684
$ws.Range("K46").Value = 2579.8
$ws.Range("L46").Value = 2132.3684
$ws.Range("M46").Value = -2391.8
$ws.Range("N46").Value = -2508.3684

$ws.Range("H61").Value = 2116
$ws.Range("I61").Value = 1863.4286
$ws.Range("K61").Value = 1863.4286
$ws.Range("M61").Value = -1661.4286

$ws.Range("H113").Value = 2116
$ws.Range("I113").Value = 1863.4286
$ws.Range("K113").Value = 1863.4286
$ws.Range("M113").Value = 306.5714

$ws.Range("H122").Value = 6526.933
$ws.Range("I122").Value = 4944.3335
$ws.Range("K122").Value = 14833.0005
$ws.Range("M122").Value = -12383.0005

$ws.Range("H123").Value = 30948
$ws.Range("J123").Value = 30948
$ws.Range("L123").Value = 30948
$ws.Range("N123").Value = -40748

$ws.Range("H125").Value = 42915
$ws.Range("J125").Value = 42915
$ws.Range("L125").Value = 42915
$ws.Range("N125").Value = -52755

$ws.Range("H126").Value = 6005.857
$ws.Range("I126").Value = 3684
$ws.Range("J126").Value = 8327.714
$ws.Range("K126").Value = 11052
$ws.Range("L126").Value = 24983.142
$ws.Range("M126").Value = -8582
$ws.Range("N126").Value = -29923.142

$ws.Range("H132").Value = 7329.3076
$ws.Range("I132").Value = 2434
$ws.Range("K132").Value = 7302
$ws.Range("M132").Value = -4772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4463.364
$ws.Range("I122").Value = 2954
$ws.Range("K122").Value = 8862
$ws.Range("M122").Value = -6412
